$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.756.55"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.734.40"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "3.218.12"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "63.601.46"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "2.737.86"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.53"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "0.0₃0904"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.32"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.09"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.06"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.974"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "346.66"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.92"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.100"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.19%  "
